# Swap the "name" identifiers that were assigned to the two logo pictures
# (Pearson PNG logo in the footers, BTEC JPG logo in the headers) that
# live in this document's default and first-page headers/footers.
#
#   Pearson logo (footers):  image2.png -> image1.png
#   BTEC logo    (headers):  image1.jpg -> image2.jpg
#
# wdHeaderFooterPrimary = 1 (the "default" header/footer)
# wdHeaderFooterFirstPage = 2 (the "first page" header/footer)

$d = $word.ActiveDocument
$section = $d.Sections.Item(1)

# Headers: BTEC_Logo-Orange (jpg) -> rename image1.jpg to image2.jpg
foreach ($idx in 1, 2) {
    $header = $section.Headers.Item($idx)
    if ($header.Exists) {
        for ($i = 1; $i -le $header.Range.InlineShapes.Count; $i++) {
            $shape = $header.Range.InlineShapes.Item($i)
            if ($shape.Name -eq "image1.jpg") {
                $shape.Name = "image2.jpg"
            }
        }
    }
}

# Footers: PearsonLogo (png) -> rename image2.png to image1.png
foreach ($idx in 1, 2) {
    $footer = $section.Footers.Item($idx)
    if ($footer.Exists) {
        for ($i = 1; $i -le $footer.Range.InlineShapes.Count; $i++) {
            $shape = $footer.Range.InlineShapes.Item($i)
            if ($shape.Name -eq "image2.png") {
                $shape.Name = "image1.png"
            }
        }
    }
}
